$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '30.386.30'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +1.29%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.925.20'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +0.86%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.9999'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '0.8084'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.25%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '244.86'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.30%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.9999'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.05%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3264'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +3.06%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '27.34'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +3.93%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.07278'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +5.75%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.8008'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +7.87%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.08098'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +1.23%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.926.06'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.94%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.425'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +4.49%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '94.64'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.79%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '30.391.13'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +1.30%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '14.39'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +3.44%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '6.126'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +4.40%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '253.38'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +3.28%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.000007883'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +1.97%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '8.103'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +18.65%  '
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("B24").Value = 'Stellar'
$ws.Range("C24").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.1635'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +17.59%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '9.581'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +3.90%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '167.52'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.37%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '19.16'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.47%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.161'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +6.23%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.377'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.87%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.547'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.93%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '4.358'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.11%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.164'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +2.04%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.05637'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +2.20%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.306'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +4.13%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.7459'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +1.73%  '
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +0.49%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.723'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01961'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +1.82%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.818'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +1.12%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.4516'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +2.22%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '74.32'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +2.94%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '6.005'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -2.15%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.942'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +3.68%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.8553'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +2.13%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.9999'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '103.72'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +3.16%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.033.96'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +4.87%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '10.02'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +2.84%  '
$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '7.678'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +1.85%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.076.18'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +1.17%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.574'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +6.57%  '
